# Generate Report for Handoff
# Refresh the localization-status report: the handoff status flips from
# "Handed back: in sync with en-US" to "Ready for handoff", the associated
# timestamps advance a few seconds, and the now-shorter status text lets the
# Status/Latest Handoff Datetime columns shrink to fit.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Timestamps bumped a few seconds forward ---
$wsOverview.Range("G2").Value = "2016-08-22 06:57:57"
$wsDeDe.Range("H2").Value     = "2016-08-22 06:57:57"
$wsZhCn.Range("H2").Value     = "2016-08-22 06:57:53"

# --- Column widths: the Status / Latest Handoff Datetime columns narrow now
#     that "Ready for handoff" is shorter than the old status text ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.333333333333332
